$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-7 and append new rows 8-13 to reflect the latest
# weekly Femacal de La Calera - Chirimoya price data.

# Row 2
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "Femacal de La Calera"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44160
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = "Otros"
$ws.Range("I2").Value = 100107002
$ws.Range("J2").Value = "Chirimoya"
$ws.Range("K2").Value = "Cultivar IV Región"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 102
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 20882
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("R2").Value = "Provincia del Elquí"
$ws.Range("S2").Value = 2088
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Femacal de La Calera"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44159
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107002
$ws.Range("J3").Value = "Chirimoya"
$ws.Range("K3").Value = "Cultivar IV Región"
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 45
$ws.Range("N3").Value = 22000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 22000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("R3").Value = "Provincia del Elquí"
$ws.Range("S3").Value = 2200
$ws.Range("T3").Value = 10

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Femacal de La Calera"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44159
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107002
$ws.Range("J4").Value = "Chirimoya"
$ws.Range("K4").Value = "Cultivar IV Región"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 47
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("R4").Value = "Provincia del Elquí"
$ws.Range("S4").Value = 2000
$ws.Range("T4").Value = 10

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Femacal de La Calera"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44441
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = "Otros"
$ws.Range("I5").Value = 100107002
$ws.Range("J5").Value = "Chirimoya"
$ws.Range("K5").Value = "Cultivar IV Región"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 68
$ws.Range("N5").Value = 3000
$ws.Range("O5").Value = 3000
$ws.Range("P5").Value = 3000
$ws.Range("Q5").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R5").Value = "Provincia del Elquí"
$ws.Range("S5").Value = 3000
$ws.Range("T5").Value = 1

# Row 6
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Femacal de La Calera"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44441
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107002
$ws.Range("J6").Value = "Chirimoya"
$ws.Range("K6").Value = "Cultivar IV Región"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 70
$ws.Range("N6").Value = 2700
$ws.Range("O6").Value = 2700
$ws.Range("P6").Value = 2700
$ws.Range("Q6").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R6").Value = "Provincia del Elquí"
$ws.Range("S6").Value = 2700
$ws.Range("T6").Value = 1

# Row 7
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "Femacal de La Calera"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 44446
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100107
$ws.Range("H7").Value = "Otros"
$ws.Range("I7").Value = 100107002
$ws.Range("J7").Value = "Chirimoya"
$ws.Range("K7").Value = "Cultivar IV Región"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 45
$ws.Range("N7").Value = 30000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 30000
$ws.Range("Q7").Value = "$/bandeja 10 kilos"
$ws.Range("R7").Value = "Provincia del Elquí"
$ws.Range("S7").Value = 3000
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "Femacal de La Calera"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44446
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = "Otros"
$ws.Range("I8").Value = 100107002
$ws.Range("J8").Value = "Chirimoya"
$ws.Range("K8").Value = "Cultivar IV Región"
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = 28000
$ws.Range("O8").Value = 28000
$ws.Range("P8").Value = 28000
$ws.Range("Q8").Value = "$/bandeja 10 kilos"
$ws.Range("R8").Value = "Provincia del Elquí"
$ws.Range("S8").Value = 2800
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = "Femacal de La Calera"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44438
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107002
$ws.Range("J9").Value = "Chirimoya"
$ws.Range("K9").Value = "Cultivar IV Región"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 85
$ws.Range("N9").Value = 27000
$ws.Range("O9").Value = 30000
$ws.Range("P9").Value = 28588
$ws.Range("Q9").Value = "$/bandeja 10 kilos"
$ws.Range("R9").Value = "Provincia del Elquí"
$ws.Range("S9").Value = 2859
$ws.Range("T9").Value = 10

# Row 10
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "Femacal de La Calera"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 44442
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100107
$ws.Range("H10").Value = "Otros"
$ws.Range("I10").Value = 100107002
$ws.Range("J10").Value = "Chirimoya"
$ws.Range("K10").Value = "Cultivar IV Región"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 3000
$ws.Range("O10").Value = 3000
$ws.Range("P10").Value = 3000
$ws.Range("Q10").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R10").Value = "Provincia del Elquí"
$ws.Range("S10").Value = 3000
$ws.Range("T10").Value = 1

# Row 11
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "Femacal de La Calera"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44442
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100107
$ws.Range("H11").Value = "Otros"
$ws.Range("I11").Value = 100107002
$ws.Range("J11").Value = "Chirimoya"
$ws.Range("K11").Value = "Cultivar IV Región"
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 473
$ws.Range("N11").Value = 2700
$ws.Range("O11").Value = 2700
$ws.Range("P11").Value = 2700
$ws.Range("Q11").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R11").Value = "Provincia del Elquí"
$ws.Range("S11").Value = 2700
$ws.Range("T11").Value = 1

# Row 12
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "Femacal de La Calera"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44435
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100107
$ws.Range("H12").Value = "Otros"
$ws.Range("I12").Value = 100107002
$ws.Range("J12").Value = "Chirimoya"
$ws.Range("K12").Value = "Cultivar IV Región"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 45
$ws.Range("N12").Value = 30000
$ws.Range("O12").Value = 30000
$ws.Range("P12").Value = 30000
$ws.Range("Q12").Value = "$/bandeja 10 kilos"
$ws.Range("R12").Value = "Provincia del Elquí"
$ws.Range("S12").Value = 3000
$ws.Range("T12").Value = 10

# Row 13
$ws.Range("A13").Value = 3
$ws.Range("B13").Value = "Femacal de La Calera"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44435
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100107
$ws.Range("H13").Value = "Otros"
$ws.Range("I13").Value = 100107002
$ws.Range("J13").Value = "Chirimoya"
$ws.Range("K13").Value = "Cultivar IV Región"
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 47
$ws.Range("N13").Value = 28000
$ws.Range("O13").Value = 28000
$ws.Range("P13").Value = 28000
$ws.Range("Q13").Value = "$/bandeja 10 kilos"
$ws.Range("R13").Value = "Provincia del Elquí"
$ws.Range("S13").Value = 2800
$ws.Range("T13").Value = 10

# Apply the date/time number format (same style used by D2:D7) to the
# date cells in the newly appended rows so they match the existing column style.
$ws.Range("D8:D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Applied weekly Fruta/hortaliza update"
